$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12, shifting rows 12-13 down to 13-14
$ws.Rows.Item(12).Insert()

# Set the new cell's content (plain formatting, default row height)
$ws.Range("B12").Value = "For test scenarios and test cases I would rather use Jira, TestRail or some other testing management tool."

# Update the selection to match the target state
$ws.Range("B12").Select()
